# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de report sheets (row 2 only) to reflect a
# freshly regenerated handback report.

$wb = $excel.ActiveWorkbook

$ws_zhcn = $wb.Worksheets.Item("zh-cn")
$ws_zhcn.Range("E2").Value2 = "2016-03-24 05:13:46"
$ws_zhcn.Range("H2").Value2 = "2016-03-24 05:14:09"

$ws_dede = $wb.Worksheets.Item("de-de")
$ws_dede.Range("E2").Value2 = "2016-03-24 05:13:51"
$ws_dede.Range("H2").Value2 = "2016-03-24 05:14:16"
